# Data cleaning pass applied before filters:
#  1. Fix header text in A1 (drop stray " 번호" suffix).
#  2. Convert the budget column (I2:I165) from comma-formatted text
#     ("483,561,215") to genuine numeric values (483561215).
#  3. Collapse a double space in the G20 service-name text.
#  4. Trim trailing whitespace from the C121 agency name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header cleanup -----------------------------------------------
$ws.Range("A1").Value = "용역 발주계획목록"

# --- 2. Column I: strip thousands separators, store as numbers -------
for ($r = 2; $r -le 165; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $text = [string]$cell.Value2
    if ($text -ne "") {
        $clean = $text -replace ",", ""
        $cell.Value = [double]$clean
    }
}

# --- 3. Whitespace / text fixes ---------------------------------------
$ws.Range("G20").Value = "경부고속도로 직선화 공사(우회도로 3, 4) 건설사업관리용역"
$ws.Range("C121").Value = "(사)대한전기협회"
